$wb = $excel.ActiveWorkbook

# The user switched to the "Repayment schedule" sheet (it becomes the
# active tab) and inserted a new blank column before column N ("Late"),
# which pushes "Late" / heading / "Outstanding" one column to the right.
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Activate()

$ws.Columns("N").Insert() | Out-Null

# Excel carries the formatting (incl. width) of the column to the left
# into a freshly inserted column.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Leave the selection where the user ended up after the insert.
$ws.Range("L12").Select() | Out-Null
